$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 176, shifting existing rows 176-210 down to 177-211
$ws.Rows.Item(176).Insert()

# Populate the new row 176 with values (copy of old row 176's constant fields + new date/price data)
$ws.Cells.Item(176, 1).Value = 7
$ws.Cells.Item(176, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(176, 3).Value = "Ñuble"
$ws.Cells.Item(176, 4).Value = 44641
$ws.Cells.Item(176, 5).Value = 16
$ws.Cells.Item(176, 6).Value = 100112003
$ws.Cells.Item(176, 7).Value = "Ajo"
$ws.Cells.Item(176, 8).Value = "Chino"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 120
$ws.Cells.Item(176, 11).Value = 20000
$ws.Cells.Item(176, 12).Value = 21000
$ws.Cells.Item(176, 13).Value = 20500
$ws.Cells.Item(176, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(176, 15).Value = "China"
$ws.Cells.Item(176, 16).Value = 2050
$ws.Cells.Item(176, 17).Value = 10
$ws.Cells.Item(176, 18).Value = "Hortaliza"
